$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 99

# Columns that look numeric/date-like must be forced to stay as plain text,
# matching the rest of the sheet's convention (e.g. A2 "-3", D2 "14" are text).
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "6394"

$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 2).Value = "7/14/2025"

$ws.Cells.Item($row, 3).Value = "LAMBARE 1076"

$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "5"

$ws.Cells.Item($row, 5).NumberFormat = "@"
$ws.Cells.Item($row, 5).Value = "808194286"

$ws.Cells.Item($row, 6).Value = "AYKO"
$ws.Cells.Item($row, 7).Value = "Pendiente"
$ws.Cells.Item($row, 8).Value = "Picada"

# Numeric columns
$ws.Cells.Item($row, 9).Value = 1

$ws.Cells.Item($row, 10).Value = "Cambio"
$ws.Cells.Item($row, 11).Value = "Sin equipos"
$ws.Cells.Item($row, 12).Value = "Pasante"

$ws.Cells.Item($row, 13).Value = -58.43008
$ws.Cells.Item($row, 14).Value = -34.601416

$ws.Cells.Item($row, 15).Value = "Almagro"
$ws.Cells.Item($row, 16).Value = "Capital Sur"
